$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.677.02"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "'1.961.48"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'248.69"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'0.4839"
$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("D8").Value = "'0.2955"
$ws.Range("E8").Value = "  +2.78%  "

$ws.Range("D9").Value = "'0.06788"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").Value = "'110.58"
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("D11").Value = "'19.34"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").Value = "'1.960.24"
$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").Value = "'0.07744"
$ws.Range("E13").Value = "  +2.39%  "

$ws.Range("D14").Value = "'5.483"
$ws.Range("E14").Value = "  +4.50%  "

$ws.Range("E15").Value = "  +3.53%  "

$ws.Range("D16").Value = "'293.23"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'30.677.74"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "'13.33"
$ws.Range("E18").Value = "  +3.35%  "

$ws.Range("D19").Value = "'5.657"
$ws.Range("E19").Value = "  +3.44%  "

$ws.Range("D20").Value = "'0.000007690"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").Value = "'2.217.77"
$ws.Range("E21").Value = "  +2.43%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'6.615"
$ws.Range("E24").Value = "  +3.52%  "

$ws.Range("D25").Value = "'9.894"
$ws.Range("E25").Value = "  +4.65%  "

$ws.Range("D26").Value = "'170.32"
$ws.Range("E26").Value = "  +3.62%  "

$ws.Range("D27").Value = "'20.15"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").Value = "'2.200"
$ws.Range("E28").Value = "  +3.35%  "

$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("D30").Value = "'1.439"
$ws.Range("E30").Value = "  +2.58%  "

$ws.Range("D31").Value = "'4.680"
$ws.Range("E31").Value = "  +16.38%  "

$ws.Range("D32").Value = "'4.461"
$ws.Range("E32").Value = "  +7.23%  "

$ws.Range("D33").Value = "'0.05119"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("D34").Value = "'0.7797"
$ws.Range("E34").Value = "  +7.02%  "

$ws.Range("D35").Value = "'1.175"
$ws.Range("E35").Value = "  +3.71%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "'2.734"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("D38").Value = "'2.711"
$ws.Range("E38").Value = "  +1.54%  "

$ws.Range("D39").Value = "'2.074"
$ws.Range("E39").Value = "  +3.44%  "

$ws.Range("D40").Value = "'111.15"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").Value = "'6.128"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("D42").Value = "'0.4461"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("D43").Value = "'0.8758"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("D44").Value = "'70.41"
$ws.Range("E44").Value = "  +3.37%  "

$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").Value = "'7.397"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").Value = "'0.1286"
$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").Value = "'9.346"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("D50").Value = "'47.75"
$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.2515"
$ws.Range("E51").Value = "  -0.17%  "
